$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.039.93"
$ws.Range("E2").Value = "  -1.02%  "
$ws.Range("D3").Value = "2.313.06"
$ws.Range("E3").Value = "  -2.02%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "532.64"
$ws.Range("E5").Value = "  +1.93%  "
$ws.Range("D6").Value = "132.16"
$ws.Range("E6").Value = "  -3.40%  "
$ws.Range("E7").Value = "  -0.17%  "
$ws.Range("E8").Value = "  -0.74%  "
$ws.Range("D9").Value = "2.334.09"
$ws.Range("E9").Value = "  -1.94%  "
$ws.Range("E10").Value = "  -1.64%  "
$ws.Range("E11").Value = "  +0.06%  "
$ws.Range("D12").Value = "5.31"
$ws.Range("E12").Value = "  -2.88%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").Value = "2.730.29"
$ws.Range("E14").Value = "  -1.91%  "
$ws.Range("D15").Value = "23.38"
$ws.Range("E15").Value = "  -4.03%  "
$ws.Range("D16").Value = "57.080.22"
$ws.Range("E16").Value = "  -0.87%  "
$ws.Range("E17").Value = "  -2.79%  "
$ws.Range("D18").Value = "2.332.24"
$ws.Range("E18").Value = "  -1.47%  "
$ws.Range("D19").Value = "337.01"
$ws.Range("E19").Value = "  +1.69%  "
$ws.Range("D20").Value = "10.41"
$ws.Range("E20").Value = "  -2.26%  "
$ws.Range("D21").Value = "6.91"
$ws.Range("E21").Value = "  +2.48%  "
$ws.Range("D22").Value = "4.15"
$ws.Range("E22").Value = "  -2.53%  "
$ws.Range("E23").Value = "  -0.04%  "
$ws.Range("D24").Value = "61.61"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("D26").Value = "8.66"
$ws.Range("E26").Value = "  +4.09%  "
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.05%  "
$ws.Range("E28").Value = "  +0.07%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "171.00"
$ws.Range("E29").Value = "  +0.35%  "
$ws.Range("D30").Value = "1.72"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  -3.35%  "
$ws.Range("D33").Value = "18.51"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("E34").Value = "  -0.03%  "
$ws.Range("D35").Value = "0.992"
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "1.25"
$ws.Range("E36").Value = "  -4.35%  "
$ws.Range("E37").Value = "  -1.71%  "
$ws.Range("E38").Value = "  -3.18%  "
$ws.Range("D39").Value = "1.58"
$ws.Range("E39").Value = "  -0.53%  "
$ws.Range("D40").Value = "39.15"
$ws.Range("E40").Value = "  +1.20%  "
$ws.Range("D41").Value = "148.53"
$ws.Range("E41").Value = "  -2.02%  "
$ws.Range("D42").Value = "0.376"
$ws.Range("E42").Value = "  -2.17%  "
$ws.Range("D43").Value = "3.59"
$ws.Range("E43").Value = "  -2.09%  "
$ws.Range("D44").Value = "280.37"
$ws.Range("E44").Value = "  -1.33%  "
$ws.Range("E45").Value = "  -5.51%  "
$ws.Range("D46").Value = "0.0926"
$ws.Range("E46").Value = "  -1.64%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0500"
$ws.Range("E47").Value = "  -1.83%  "
$ws.Range("D48").Value = "18.69"
$ws.Range("E48").Value = "  +1.41%  "
$ws.Range("E49").Value = "  -1.96%  "
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("E51").Value = "  +5.50%  "
